$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "efficiency" column header for the two performance tables ---
$ws.Range("C13").Value = "efficiency"
$ws.Range("C13").Font.Color = 255

$ws.Range("C26").Value = "efficiency"
$ws.Range("C26").Font.Color = 255

$ws.Range("D14").Value = 1
$ws.Range("D15").Value = 2
$ws.Range("D16").Value = 4
$ws.Range("D17").Value = 8
$ws.Range("D18").Value = 1
$ws.Range("D19").Value = 2
$ws.Range("D20").Value = 4
$ws.Range("D21").Value = 8
$ws.Range("C14").Formula = '=$B$14/(B14*D14)'
$ws.Range("C14").Font.Color = 255
$ws.Range("C15").Formula = '=$B$14/(B15*D15)'
$ws.Range("C15").Font.Color = 255
$ws.Range("C16").Formula = '=$B$14/(B16*D16)'
$ws.Range("C16").Font.Color = 255
$ws.Range("C17").Formula = '=$B$14/(B17*D17)'
$ws.Range("C17").Font.Color = 255
$ws.Range("C18").Formula = '=$B$14/(B18*D18)'
$ws.Range("C18").Font.Color = 255
$ws.Range("C19").Formula = '=$B$14/(B19*D19)'
$ws.Range("C19").Font.Color = 255
$ws.Range("C20").Formula = '=$B$14/(B20*D20)'
$ws.Range("C20").Font.Color = 255
$ws.Range("C21").Formula = '=$B$14/(B21*D21)'
$ws.Range("C21").Font.Color = 255
$ws.Range("D27").Value = 1
$ws.Range("D28").Value = 2
$ws.Range("D29").Value = 4
$ws.Range("D30").Value = 8
$ws.Range("D31").Value = 1
$ws.Range("D32").Value = 2
$ws.Range("D33").Value = 4
$ws.Range("D34").Value = 8
$ws.Range("D35").Value = 1
$ws.Range("D36").Value = 2
$ws.Range("D37").Value = 4
$ws.Range("D38").Value = 8
$ws.Range("D39").Value = 1
$ws.Range("D40").Value = 2
$ws.Range("D41").Value = 4
$ws.Range("D42").Value = 8
$ws.Range("D43").Value = 1
$ws.Range("D44").Value = 2
$ws.Range("D45").Value = 4
$ws.Range("D46").Value = 8
$ws.Range("D47").Value = 1
$ws.Range("D48").Value = 2
$ws.Range("D49").Value = 4
$ws.Range("D50").Value = 8
$ws.Range("C27").Formula = '=$B$27/(B27*D27)'
$ws.Range("C27").Font.Color = 255
$ws.Range("C28").Formula = '=$B$27/(B28*D28)'
$ws.Range("C28").Font.Color = 255
$ws.Range("C29").Formula = '=$B$27/(B29*D29)'
$ws.Range("C29").Font.Color = 255
$ws.Range("C30").Formula = '=$B$27/(B30*D30)'
$ws.Range("C30").Font.Color = 255
$ws.Range("C31").Formula = '=$B$27/(B31*D31)'
$ws.Range("C31").Font.Color = 255
$ws.Range("C32").Formula = '=$B$27/(B32*D32)'
$ws.Range("C32").Font.Color = 255
$ws.Range("C33").Formula = '=$B$27/(B33*D33)'
$ws.Range("C33").Font.Color = 255
$ws.Range("C34").Formula = '=$B$27/(B34*D34)'
$ws.Range("C34").Font.Color = 255
$ws.Range("C35").Formula = '=$B$27/(B35*D35)'
$ws.Range("C35").Font.Color = 255
$ws.Range("C36").Formula = '=$B$27/(B36*D36)'
$ws.Range("C36").Font.Color = 255
$ws.Range("C37").Formula = '=$B$27/(B37*D37)'
$ws.Range("C37").Font.Color = 255
$ws.Range("C38").Formula = '=$B$27/(B38*D38)'
$ws.Range("C38").Font.Color = 255
$ws.Range("C39").Formula = '=$B$27/(B39*D39)'
$ws.Range("C39").Font.Color = 255
$ws.Range("C40").Formula = '=$B$27/(B40*D40)'
$ws.Range("C40").Font.Color = 255
$ws.Range("C41").Formula = '=$B$27/(B41*D41)'
$ws.Range("C41").Font.Color = 255
$ws.Range("C42").Formula = '=$B$27/(B42*D42)'
$ws.Range("C42").Font.Color = 255
$ws.Range("C43").Formula = '=$B$27/(B43*D43)'
$ws.Range("C43").Font.Color = 255
$ws.Range("C44").Formula = '=$B$27/(B44*D44)'
$ws.Range("C44").Font.Color = 255
$ws.Range("C45").Formula = '=$B$27/(B45*D45)'
$ws.Range("C45").Font.Color = 255
$ws.Range("C46").Formula = '=$B$27/(B46*D46)'
$ws.Range("C46").Font.Color = 255
$ws.Range("C47").Formula = '=$B$27/(B47*D47)'
$ws.Range("C47").Font.Color = 255
$ws.Range("C48").Formula = '=$B$27/(B48*D48)'
$ws.Range("C48").Font.Color = 255
$ws.Range("C49").Formula = '=$B$27/(B49*D49)'
$ws.Range("C49").Font.Color = 255
$ws.Range("C50").Formula = '=$B$27/(B50*D50)'
$ws.Range("C50").Font.Color = 255
